# 242: myr changes
# - Rename "Section 17(4)" -> "Vehicle Statistics" and give it a header row
#   describing the vehicle statistics columns that the new sheet will hold.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Section 17(4)")
$ws.Name = "Vehicle Statistics"

# Header values - set in this particular order so that new shared-string
# entries are created in the same order as the target workbook
# (Make, Model Name, Range, ZEV Type, Submitted Count, Issued Count).
$ws.Range("A1").Value = "Vehicle Class"
$ws.Range("B1").Value = "ZEV Class"
$ws.Range("C1").Value = "Make"
$ws.Range("D1").Value = "Model Name"
$ws.Range("E1").Value = "Model Year"
$ws.Range("G1").Value = "Range"
$ws.Range("F1").Value = "ZEV Type"
$ws.Range("H1").Value = "Submitted Count"
$ws.Range("I1").Value = "Issued Count"

# Header row is bold, matching the look of the other sheets in the workbook.
$ws.Range("A1:I1").Font.Bold = $true

# Column widths matching the template (values chosen so the saved XML
# "width" comes out as close as possible to the source template's values).
$ws.Columns.Item(1).ColumnWidth = 15.333333333333334
$ws.Columns.Item(2).ColumnWidth = 15.333333333333334
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 17.666666666666668
$ws.Columns.Item(5).ColumnWidth = 16.166666666666668
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(7).ColumnWidth = 19.333333333333332
$ws.Columns.Item(8).ColumnWidth = 22.0
$ws.Columns.Item(9).ColumnWidth = 24.0
